$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grouped matches")

$ws.Range("B2").Value = "P8260"
$ws.Range("C2").Value = "{'eft:dpal-dbyangs'}"
$ws.Range("B3").Value = "P8266"
$ws.Range("C3").Value = "{'eft:dharmatasila', 'eft:ch-nyi-tsultrim'}"
$ws.Range("B4").Value = "P2956"
$ws.Range("C4").Value = "{'eft:krsnapandita'}"
$ws.Range("B5").Value = "https://lod.dila.edu.tw/resource.php?id=A000089"
$ws.Range("C5").Value = "{'eft:siladharma'}"
$ws.Range("B6").Value = "P8205"
$ws.Range("C6").Value = "{'eft:band-yesh-de', 'eft:yesh-d-', 'eft:zhang-yesh-d-', 'eft:ye-shes-sde', 'eft:band-yesh-d-', 'eft:yesh-d-ye-shes-sde-'}"
$ws.Range("B7").Value = "P4242"
$ws.Range("C7").Value = "{'eft:sherab-lekpa'}"
$ws.Range("B8").Value = "P3214"
$ws.Range("C8").Value = "{'eft:danasila'}"
$ws.Range("B9").Value = "P4259"
$ws.Range("C9").Value = "{'eft:dpal-gyi-lhun-po', 'eft:ban-de-dpal-gyi-lhun-po', 'eft:palgyi-lh-npo'}"
$ws.Range("B10").Value = "?"
$ws.Range("C10").Value = "{'eft:sakyasena'}"
$ws.Range("B11").Value = "P0TMP092"
$ws.Range("C11").Value = "{'eft:anandasri-s-'}"
$ws.Range("B12").Value = "P4CZ16780"
$ws.Range("C12").Value = "{'eft:manjusrigarbha'}"
$ws.Range("B13").Value = "P8182"
$ws.Range("C13").Value = "{'eft:dpal-brtsegs', 'eft:ban-de-dpal-brtsegs', 'eft:paltsek', 'eft:ska-ba-dpal-brtsegs', 'eft:kawa-paltsek-under-the-name-paltsek-raksita-'}"
$ws.Range("B14").Value = "P2548"
$ws.Range("C14").Value = "{'eft:prajnavarma', 'eft:prajnavarman'}"
$ws.Range("B15").Value = "P8222"
$ws.Range("C15").Value = "{'eft:jnanasiddhi', 'eft:jnanasidhi'}"
$ws.Range("B16").Value = "P3709"
$ws.Range("C16").Value = "{'eft:phakpa-sherab'}"
$ws.Range("B17").Value = "P0TMP080"
$ws.Range("C17").Value = "{'eft:hwa-shang-zab-mo'}"
$ws.Range("B18").Value = "P8261"
$ws.Range("C18").Value = "{'eft:munivarma', 'eft:munivarman'}"
$ws.Range("B19").Value = "P8249"
$ws.Range("C19").Value = "{'eft:dharmakara'}"
$ws.Range("B20").Value = "P3379"
$ws.Range("C20").Value = "{'eft:dipamkarasrijnana', 'eft:dipamkara-srijnana'}"
$ws.Range("B21").Value = "P8213"
$ws.Range("C21").Value = "{'eft:t-vidyakarasimha', 'eft:vidyakarasimha'}"
$ws.Range("B22").Value = "P753"
$ws.Range("C22").Value = "{'eft:rin-chen-bzang-po'}"
$ws.Range("B23").Value = "P4258"
$ws.Range("C23").Value = "{'eft:dpal-byor'}"
$ws.Range("B24").Value = "P0TMP104"
$ws.Range("C24").Value = "{'eft:punyasambhava'}"
$ws.Range("B25").Value = "P0RK8"
$ws.Range("C25").Value = "{'eft:dharmapala'}"
$ws.Range("B26").Value = "P8211"
$ws.Range("C26").Value = "{'eft:vidyakaraprabha'}"
$ws.Range("B27").Value = "P8171"
$ws.Range("C27").Value = "{'eft:dharmasribhadra'}"
$ws.Range("B28").Value = "P5651"
$ws.Range("C28").Value = "{'eft:pa-tshab-nyi-ma-grags'}"
$ws.Range("B29").Value = "P4255"
$ws.Range("C29").Value = "{'eft:yesh-nyingpo', 'eft:ye-shes-snying-po', 'eft:t-jnanagarbha'}"
$ws.Range("B30").Value = "P4263"
$ws.Range("C30").Value = "{'eft:dge-ba-dpal'}"
$ws.Range("B31").Value = "P8219"
$ws.Range("C31").Value = "{'eft:visuddhasimha'}"
$ws.Range("B32").Value = "P8217"
$ws.Range("C32").Value = "{'eft:jnanagarbha', 'eft:t-jnanagarbha'}"
$ws.Range("B33").Value = "P8269"
$ws.Range("C33").Value = "{'eft:dgon-gling-rma'}"
$ws.Range("B34").Value = "P8265"
$ws.Range("C34").Value = "{'eft:ratnaraksita'}"
$ws.Range("B35").Value = "P2637"
$ws.Range("C35").Value = "{'eft:trakpa-gyaltsen'}"
$ws.Range("B36").Value = "P8267"
$ws.Range("C36").Value = "{'eft:vijayasila'}"
$ws.Range("B37").Value = "P3456"
$ws.Range("C37").Value = "{'eft:tshul-khrims-rgyal-ba'}"
$ws.Range("B38").Value = "P8183"
$ws.Range("C38").Value = "{'eft:cog-ro-klu-i-rgyal-mtshan', 'eft:klu-i-rgyal-mtshan'}"
$ws.Range("B39").Value = "P4CZ15137"
$ws.Range("C39").Value = "{'eft:kumarakalasa'}"
$ws.Range("B40").Value = "P8209"
$ws.Range("C40").Value = "{'eft:jinamitra-k-', 'eft:dzi-na-mi-tra-k-', 'eft:jinamitra'}"
$ws.Range("B41").Value = "P3285"
$ws.Range("C41").Value = "{'eft:sakya-yesh-'}"
$ws.Range("B42").Value = "P8245"
$ws.Range("C42").Value = "{'eft:buddhakaravarma'}"
$ws.Range("B43").Value = "P1KG8854"
$ws.Range("C43").Value = "{'eft:surendrabodhi', 'eft:srilendrabodhi', 'eft:silendrabodhi'}"
$ws.Range("B44").Value = "P8210"
$ws.Range("C44").Value = "{'eft:danasila'}"
$ws.Range("B45").Value = "P8220"
$ws.Range("C45").Value = "{'eft:devacandra'}"
$ws.Range("B46").Value = "P8151"
$ws.Range("C46").Value = "{'eft:gayadhara'}"
$ws.Range("B47").Value = "P4CZ16819"
$ws.Range("C47").Value = "{'eft:sakyaprabha'}"
$ws.Range("B48").Value = "P00KG07267"
$ws.Range("C48").Value = "{'eft:sarvanyadeva', 'eft:sarvajnadeva'}"
$ws.Range("B49").Value = "P8228"
$ws.Range("C49").Value = "{'eft:surendrabodhi'}"
$ws.Range("B50").Value = "P0TMP098"
$ws.Range("C50").Value = "{'eft:jinavara'}"
$ws.Range("B51").Value = "P8268"
$ws.Range("C51").Value = "{'eft:buddhaprabha'}"
$ws.Range("B52").Value = "P8273"
$ws.Range("C52").Value = "{'eft:rinchen-tso', 'eft:rin-chen-tsho'}"
$ws.Range("B53").Value = "P0TMPT007"
$ws.Range("C53").Value = "{'eft:rnam-par-mi-rtog-pa'}"
$ws.Range("B54").Value = "P8206"
$ws.Range("C54").Value = "{'eft:celu'}"
$ws.Range("B55").Value = "P8263"
$ws.Range("C55").Value = "{'eft:leki-d-'}"
$ws.Range("B56").Value = "P8093"
$ws.Range("C56").Value = "{'eft:kamalagupta'}"
